# Insert a new weekly record for "Haba" (Vega Modelo de Temuco) as row 74,
# pushing the existing rows 74:80 down to 75:81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("74:74").Insert()

$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44858
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100112026
$ws.Range("G74").Value = "Haba"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 155
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = 10000
$ws.Range("N74").Value = "`$/saco 25 kilos"
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 400
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = "Hortaliza"
